$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) Create the two brand-new sheets first, by duplicating a sheet that
#    already has the right sheetPr/sheetFormatPr/drawing wiring (so the
#    new sheets don't end up with bare Excel defaults). "works-section"
#    is blank and has no custom <cols>, so it's a perfect template.
# -----------------------------------------------------------------------
$template = $wb.Worksheets.Item("works-section")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "contact-section-tmp"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet2)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "admins-tmp"

# -----------------------------------------------------------------------
# 2) Rename sheets (carefully ordered so we never collide with a name
#    that's still in use) so the final tab order/names are:
#    home-section, about-section, projects-section, collaborators-section,
#    contact-section, admins
# -----------------------------------------------------------------------
$wb.Worksheets.Item("about-section").Name = "projects-section"
$wb.Worksheets.Item("admins").Name = "about-section"
$wb.Worksheets.Item("works-section").Name = "collaborators-section"
$wb.Worksheets.Item("contact-section-tmp").Name = "contact-section"
$wb.Worksheets.Item("admins-tmp").Name = "admins"

$home = $wb.Worksheets.Item("home-section")
$about = $wb.Worksheets.Item("about-section")
$projects = $wb.Worksheets.Item("projects-section")
$collaborators = $wb.Worksheets.Item("collaborators-section")
$contact = $wb.Worksheets.Item("contact-section")
$admins = $wb.Worksheets.Item("admins")

# -----------------------------------------------------------------------
# 3) home-section: drop the C-column hyperlink / url row, add a
#    "background" header, new sample text, narrower C column.
# -----------------------------------------------------------------------
foreach ($hl in $home.Hyperlinks) { $hl.Delete() }

$home.Columns.Item(3).ColumnWidth = 15.917

$home.Range("A1").Value = "textblock"
$home.Range("B1").Value = "background"
$home.Range("C1").ClearContents()

$home.Range("A2").Value = "sometext example"
$home.Range("B2").Value = "yes"
$home.Range("C2").ClearContents()

# -----------------------------------------------------------------------
# 4) about-section: heading / textblock / background row layout
# -----------------------------------------------------------------------
$home.Range("A1").Copy()
$about.Range("A1:C2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$about.Columns.Item(2).ColumnWidth = 59.084

$about.Range("A1").Value = "heading"
$about.Range("B1").Value = "textblock"
$about.Range("C1").Value = "background"

$about.Range("A2").Value = "About Us"
$about.Range("B2").Value = 'Lunar Studio is an alliance between skilled artists, writers and thinkers, to develop and publish creative artifacts. /n The Studio is an imagined space, carved out of busy, pressuredlives. It is intended to incubate wonderful, strange ideas and turn them into items which can be shared with others. /n Our intention is to enable our own independent creative work,and to work with people we really enjoy. We undertake passion projects without commercial backing by producing, funding and publishing collaboratively. /n We don''t cater to anyones'' needs or desires but our own and each others'' in the process of creating. The products which emerge are allowed to be uncompromising, decadent, and a pleasure to make.'
$about.Range("C2").Value = "yes"

# -----------------------------------------------------------------------
# 5) projects-section: id/order + name table
# -----------------------------------------------------------------------
$home.Range("A1").Copy()
$projects.Range("A1:B4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$projects.Columns.Item(2).ColumnWidth = 26.584

$projects.Range("A1").Value = "id/order"
$projects.Range("B1").Value = "name"

$projects.Range("A2").Value = 1
$projects.Range("B2").Value = "Frankie and the moon"

$projects.Range("A3").Value = 2
$projects.Range("B3").Value = "Prints and Concept Art"

$projects.Range("A4").Value = 3
$projects.Range("B4").Value = "Bone Cards"

# -----------------------------------------------------------------------
# 6) collaborators-section: requires / collaborators labels
# -----------------------------------------------------------------------
$home.Range("A1").Copy()
$collaborators.Range("A1:A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$collaborators.Range("A1").Value = "requires"
$collaborators.Range("A2").Value = "collaborators"

# -----------------------------------------------------------------------
# 7) contact-section: heading / Contact Us
# -----------------------------------------------------------------------
$home.Range("A1").Copy()
$contact.Range("A1:A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$contact.Range("A1").Value = "heading"
$contact.Range("A2").Value = "Contact Us"

# admins stays empty (new, blank sheet) - nothing further to do.
